# Updated capital structure database
# Apply updated values to rows 2 and 3 of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D"  = -0.0593
    "G"  = -0.4208333333333333
    "H"  = -0.4208333333333333
    "I"  = -0.4979166666666666
    "J"  = -0.4979166666666666
    "K"  = -26
    "L"  = -0.5416666666666666
    "U"  = 8.02
    "V"  = 1.147353361945636
    "W"  = -12.32227488151659
    "X"  = 0.7665345783699571
    "Y"  = -13.08880945988654
    "Z"  = 0.9914692334703488
    "AA" = -0.4936690558321111
    "AB" = 0.1038806531387532
    "AC" = -0.5975497089708642
    "AD" = 75.7
    "AF" = 75.7
    "AG" = 67.68000000000001
    "AH" = 0.9154674083927924
    "AI" = 1.490157480314961
    "AJ" = 0.9063881076737647
    "AK" = 1.582047685834502
    "AL" = 2.72
    "AM" = 2.697
    "AN" = -3.604761904761905
    "AO" = -8.786764705882351
    "AP" = -3.222857142857143
    "AQ" = -8.86169818316648
}

foreach ($row in 2, 3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
